$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'313.16"
$ws.Range("E2").Value = "'2.32%"
$ws.Range("D3").Value = "'37.62"
$ws.Range("E3").Value = "'0.93%"
$ws.Range("D4").Value = "'5.128"
$ws.Range("E4").Value = "'0.44%"
$ws.Range("D5").Value = "'0.07908"
$ws.Range("E5").Value = "'1.98%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.413"
$ws.Range("E6").Value = "'0.43%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.932"
$ws.Range("E7").Value = "'1.02%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.273"
$ws.Range("E8").Value = "'0.59%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.921"
$ws.Range("E9").Value = "'-7.32%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9211"
$ws.Range("E10").Value = "'-0.07%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1232"
$ws.Range("E11").Value = "'-3.22%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1926"
$ws.Range("E12").Value = "'1.70%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09144"
$ws.Range("E13").Value = "'3.96%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03313"
$ws.Range("E14").Value = "'-3.88%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09600"
$ws.Range("E15").Value = "'-1.50%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001389"
$ws.Range("E16").Value = "'1.41%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005795"
$ws.Range("E17").Value = "'-5.31%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.512"
$ws.Range("E18").Value = "'-1.45%"
$ws.Range("D19").Value = "'0.3444"
$ws.Range("E19").Value = "'2.09%"
$ws.Range("D20").Value = "'5.260"
$ws.Range("E20").Value = "'4.28%"
$ws.Range("E21").Value = "'-1.17%"
$ws.Range("D22").Value = "'0.2590"
$ws.Range("E22").Value = "'3.52%"
$ws.Range("E23").Value = "'-0.70%"
$ws.Range("D24").Value = "'0.04366"
$ws.Range("E24").Value = "'0.47%"
$ws.Range("E25").Value = "'1.82%"
$ws.Range("D26").Value = "'0.004303"
$ws.Range("E26").Value = "'-4.24%"
$ws.Range("E27").Value = "'-10.25%"
$ws.Range("D39").Value = "'0.02254"
$ws.Range("E39").Value = "'3.80%"
$ws.Range("D40").Value = "'0.05120"
$ws.Range("E40").Value = "'3.67%"
$ws.Range("D41").Value = "'0.007451"
$ws.Range("E41").Value = "'-3.38%"
$ws.Range("D42").Value = "'0.1363"
$ws.Range("E42").Value = "'1.65%"
$ws.Range("D43").Value = "'0.008785"
$ws.Range("E43").Value = "'-10.83%"
$ws.Range("E44").Value = "'0.09%"
$ws.Range("D45").Value = "'0.008611"
$ws.Range("E45").Value = "'-2.89%"
$ws.Range("D46").Value = "'0.00006735"
$ws.Range("E46").Value = "'-1.55%"
$ws.Range("E47").Value = "'-0.72%"
$ws.Range("D48").Value = "'0.003350"
$ws.Range("E48").Value = "'10.95%"
$ws.Range("D49").Value = "'0.001200"
$ws.Range("E49").Value = "'-8.35%"
$ws.Range("E50").Value = "'-0.72%"
$ws.Range("E51").Value = "'-0.72%"
